$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 139.25
$ws.Range("J2").Value = 200
$ws.Range("L2").Value = 200
$ws.Range("N2").Value = -426

# Row 111
$ws.Range("H111").Value = 1236.8572
$ws.Range("I111").Value = 1031.6
$ws.Range("J111").Value = 1750
$ws.Range("K111").Value = 3094.8
$ws.Range("L111").Value = 5250
$ws.Range("M111").Value = -27.79999999999973
$ws.Range("N111").Value = -11384

# Row 125
$ws.Range("H125").Value = 1610
$ws.Range("I125").Value = 2020
$ws.Range("K125").Value = 18180
$ws.Range("M125").Value = -15720

$ws = $wb.Worksheets.Item("ARM")
# Row 6
$ws.Range("H6").Value = 980
$ws.Range("I6").Value = 980
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 980
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = ""
$ws.Range("N6").Value = -807

# Row 10
$ws.Range("H10").Value = 452
$ws.Range("I10").Value = 452
$ws.Range("K10").Value = 452
$ws.Range("M10").Value = -282

# Row 29
$ws.Range("H29").Value = 10000
$ws.Range("J29").Value = 10000
$ws.Range("L29").Value = 10000
$ws.Range("N29").Value = -10616

# Row 61
$ws.Range("H61").Value = 3038.6274
$ws.Range("I61").Value = 1962.2667
$ws.Range("K61").Value = 1962.2667
$ws.Range("M61").Value = -1750.2667

# Row 63
$ws.Range("H63").Value = 4320.871
$ws.Range("I63").Value = 4666.28
$ws.Range("J63").Value = 2881.6667
$ws.Range("K63").Value = 4666.28
$ws.Range("L63").Value = 2881.6667
$ws.Range("M63").Value = -3980.28
$ws.Range("N63").Value = -4253.6667

# Row 66
$ws.Range("H66").Value = 4320.871
$ws.Range("I66").Value = 4666.28
$ws.Range("J66").Value = 2881.6667
$ws.Range("K66").Value = 23331.4
$ws.Range("L66").Value = 14408.3335
$ws.Range("M66").Value = -19899.4
$ws.Range("N66").Value = -21272.3335

# Row 132
$ws.Range("H132").Value = 3830.0977
$ws.Range("I132").Value = 4247.0435
$ws.Range("J132").Value = 3297.3333
$ws.Range("K132").Value = 12741.1305
$ws.Range("L132").Value = 9891.999899999999
$ws.Range("M132").Value = -10211.1305
$ws.Range("N132").Value = -14951.9999

# Row 136
$ws.Range("H136").Value = 3038.6274
$ws.Range("I136").Value = 1962.2667
$ws.Range("K136").Value = 5886.800099999999
$ws.Range("M136").Value = -3336.800099999999

$ws = $wb.Worksheets.Item("BSM")
# Row 133
$ws.Range("H133").Value = 30256.666
$ws.Range("J133").Value = 30256.666
$ws.Range("L133").Value = 30256.666
$ws.Range("N133").Value = -40376.666

# Row 134
$ws.Range("H134").Value = 67683.766
$ws.Range("I134").Value = 175570.67
$ws.Range("J134").Value = 8836.362999999999
$ws.Range("K134").Value = 526712.01
$ws.Range("L134").Value = 26509.089
$ws.Range("M134").Value = -524177.01
$ws.Range("N134").Value = -31579.089

$ws = $wb.Worksheets.Item("CRP")
# Row 10
$ws.Range("H10").Value = 771.2
$ws.Range("I10").Value = 714.25
$ws.Range("J10").Value = 999
$ws.Range("K10").Value = 714.25
$ws.Range("L10").Value = 999
$ws.Range("M10").Value = -575.25
$ws.Range("N10").Value = -1277

# Row 12
$ws.Range("H12").Value = 1662.2
$ws.Range("I12").Value = 768.3333
$ws.Range("J12").Value = 3003
$ws.Range("K12").Value = 768.3333
$ws.Range("L12").Value = 3003
$ws.Range("M12").Value = -598.3333
$ws.Range("N12").Value = -3343

# Row 133
$ws.Range("H133").Value = 28762
$ws.Range("J133").Value = 32495
$ws.Range("L133").Value = 32495
$ws.Range("N133").Value = -37555

# Row 134
$ws.Range("H134").Value = 2050
$ws.Range("I134").Value = 1207.2667
$ws.Range("J134").Value = 2793.5881
$ws.Range("K134").Value = 3621.800099999999
$ws.Range("L134").Value = 8380.764299999999
$ws.Range("M134").Value = -1086.800099999999
$ws.Range("N134").Value = -13450.7643

$ws = $wb.Worksheets.Item("CUL")
# Row 13
$ws.Range("H13").Value = 132.66667
$ws.Range("I13").Value = 99
$ws.Range("J13").Value = 200
$ws.Range("K13").Value = 297
$ws.Range("L13").Value = 600
$ws.Range("M13").Value = -129
$ws.Range("N13").Value = -936

# Row 56
$ws.Range("H56").Value = 6450
$ws.Range("I56").Value = 6450
$ws.Range("K56").Value = 6450
$ws.Range("M56").Value = -5920

# Row 70
$ws.Range("H70").Value = 4711.2
$ws.Range("I70").Value = 906
$ws.Range("J70").Value = 5662.5
$ws.Range("K70").Value = 2718
$ws.Range("L70").Value = 16987.5
$ws.Range("M70").Value = -2403
$ws.Range("N70").Value = -17617.5

# Row 73
$ws.Range("H73").Value = 4711.2
$ws.Range("I73").Value = 906
$ws.Range("J73").Value = 5662.5
$ws.Range("K73").Value = 2718
$ws.Range("L73").Value = 16987.5
$ws.Range("M73").Value = -1626
$ws.Range("N73").Value = -19171.5

# Row 75
$ws.Range("H75").Value = 1758.75
$ws.Range("J75").Value = 1938.5714
$ws.Range("L75").Value = 5815.7142
$ws.Range("N75").Value = -7811.7142

# Row 78
$ws.Range("H78").Value = 1758.75
$ws.Range("J78").Value = 1938.5714
$ws.Range("L78").Value = 17447.1426
$ws.Range("N78").Value = -27431.1426

# Row 131
$ws.Range("H131").Value = 2245.5454
$ws.Range("J131").Value = 1808.5
$ws.Range("L131").Value = 5425.5
$ws.Range("N131").Value = -15505.5

$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 362.14285
$ws.Range("I3").Value = 467.66666
$ws.Range("J3").Value = 283
$ws.Range("K3").Value = 467.66666
$ws.Range("L3").Value = 283
$ws.Range("M3").Value = -351.66666
$ws.Range("N3").Value = -515

$ws = $wb.Worksheets.Item("LTW")
# Row 12
$ws.Range("H12").Value = 1003
$ws.Range("I12").Value = 1003
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 1003
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = ""
$ws.Range("N12").Value = -833

# Row 129
$ws.Range("H129").Value = 23665
$ws.Range("J129").Value = 23665
$ws.Range("L129").Value = 23665
$ws.Range("N129").Value = -33665

$ws = $wb.Worksheets.Item("WVR")
# Row 33
$ws.Range("H33").Value = 4120
$ws.Range("I33").Value = 3600
$ws.Range("J33").Value = 4250
$ws.Range("K33").Value = 3600
$ws.Range("L33").Value = 4250
$ws.Range("M33").Value = -3350
$ws.Range("N33").Value = -4750

# Row 36
$ws.Range("H36").Value = 4120
$ws.Range("I36").Value = 3600
$ws.Range("J36").Value = 4250
$ws.Range("K36").Value = 3600
$ws.Range("L36").Value = 4250
$ws.Range("M36").Value = -3350
$ws.Range("N36").Value = -4750

# Row 40
$ws.Range("H40").Value = 9357.143
$ws.Range("J40").Value = 9357.143
$ws.Range("L40").Value = 9357.143
$ws.Range("N40").Value = -9655.143

# Row 132
$ws.Range("H132").Value = 19041.5
$ws.Range("I132").Value = 31371.092
$ws.Range("J132").Value = 2766.44
$ws.Range("K132").Value = 94113.276
$ws.Range("L132").Value = 8299.32
$ws.Range("M132").Value = -91583.276
$ws.Range("N132").Value = -13359.32
